$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Validación" column (column C) entirely, shifting the
# remaining columns (D:L) one position to the left.
$ws.Columns("C:C").Delete()
